$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name / surname / CvLAC URL values (fixing spanish-character issue)
$ws.Range("A2").Value = 123
$ws.Range("B2").Value = "Manuel"
$ws.Range("C2").Value = "Neira Embus"
$ws.Range("D2").Value = "http://scienti.colciencias.gov.co:8081/cvlac/visualizador/generarCurriculoCv.do?cod_rh=0001545295"

# Apply a (new) style to D2 - same font/fill/border as default, just with applyFont set
$ws.Range("D2").Font.Name = $ws.Range("D2").Font.Name

# Move selection / top-left cell to reflect where the user ended up looking
$ws.Range("D2").Select()
$excel.ActiveWindow.ScrollColumn = 2
